$wb = $excel.ActiveWorkbook

# --- Rename sheets (task order IDs refreshed) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16502912666328957"
$wb.Worksheets.Item(2).Name = "NB_TO-16502912683587701"
$wb.Worksheets.Item(3).Name = "RS_TO-16502912683607008"
$wb.Worksheets.Item(4).Name = "TOL_TO-16502912684132128"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16502912684746692"

# --- Sheet 1: GNG ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16502912665885005.csv"
$ws1.Range("B3").Value = "GNG_stims-16502912666080482.csv"
$ws1.Range("B4").Value = "go_stims-1650291266609117.csv"
$ws1.Range("B5").Value = "GNG_stims-1650291266631895.csv"

# --- Sheet 2: NB ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-16502912681978135.csv"
$ws2.Range("B3").Value = "OB-16502912676037343.csv"
$ws2.Range("B4").Value = "ZB-match_2-16502912669404364.csv"
$ws2.Range("B5").Value = "ZB-match_4-16502912667027833.csv"
$ws2.Range("B6").Value = "OB-16502912680881257.csv"
$ws2.Range("B7").Value = "TB-16502912681246765.csv"
$ws2.Range("B8").Value = "OB-16502912678997247.csv"
$ws2.Range("B9").Value = "ZB-match_4-16502912671291032.csv"
$ws2.Range("B10").Value = "TB-16502912683378227.csv"

# --- Sheet 3: RS ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4: TOL ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16502912683741748.csv"
$ws4.Range("B3").Value = "ZM_stims-1650291268363691.csv"
$ws4.Range("B4").Value = "MM_stims-1650291268396295.csv"
$ws4.Range("B5").Value = "ZM_stims-16502912683751812.csv"
$ws4.Range("B6").Value = "MM_stims-16502912684122114.csv"
$ws4.Range("B7").Value = "ZM_stims-16502912683972573.csv"

# --- Sheet 5: vSAT ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16502912684587746.csv"
$ws5.Range("B3").Value = "SAT_stims-16502912684163566.csv"
$ws5.Range("B4").Value = "vSAT_stims-1650291268443869.csv"
$ws5.Range("B5").Value = "SAT_stims-1650291268428232.csv"
